$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.324.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.194.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.193.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  +4.66%  "
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.718.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.298.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.193.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +14.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.04%  "
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "491.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0892"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.90%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.61%  "
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0659"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.915.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +12.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.52%  "

Write-Output "Applied 89 cell updates"
